# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" (strikeout) values replacing the old Strike# values in column G
$kValues = @{
    2  = 2
    3  = 1
    4  = 1
    5  = 3
    6  = 2
    7  = 1
    8  = 0
    9  = 1
    10 = 1
    11 = 2
    12 = 1
    13 = 1
    14 = 3
    15 = 4
    16 = 1
    17 = 4
    18 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
